# Updates cryptos list price/volume(1h) figures (and the Fetch.AI /
# EthereumClassic row order swap) to match the latest scraped snapshot.
# Values that are pure numbers are prefixed with a leading apostrophe so
# Excel keeps them as text (matching the original inlineStr cell layout)
# instead of silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.313.29"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "3.467.15"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "'147.62"
$ws.Range("E6").Value = "  +2.07%  "
$ws.Range("D7").Value = "3.466.68"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("E12").Value = "  +5.00%  "
$ws.Range("D13").Value = "4.062.17"
$ws.Range("D14").Value = "'29.46"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("E15").Value = "  +2.47%  "
$ws.Range("D16").Value = "3.463.98"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "63.293.47"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("E19").Value = "  +3.67%  "
$ws.Range("D20").Value = "'14.52"
$ws.Range("E20").Value = "  +3.37%  "
$ws.Range("D21").Value = "'9.33"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").Value = "'388.93"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").Value = "'0.566"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "3.612.82"
$ws.Range("D27").Value = "'0.0000116"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'8.21"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.34"
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'23.47"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D37").Value = "'7.16"
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("E38").Value = "  +4.90%  "
$ws.Range("D39").Value = "'32.07"
$ws.Range("E39").Value = "  +9.92%  "
$ws.Range("D40").Value = "'167.68"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "3.505.74"
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("D42").Value = "'0.0776"
$ws.Range("E42").Value = "  +2.83%  "
$ws.Range("D43").Value = "'0.793"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("E44").Value = "  +4.51%  "
$ws.Range("D45").Value = "'42.40"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("E46").Value = "  +3.54%  "
$ws.Range("D47").Value = "'4.38"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").Value = "2.590.22"
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("E49").Value = "  +8.94%  "
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("D51").Value = "'23.06"
$ws.Range("E51").Value = "  +0.39%  "
